# "Generate Report for Archive"
# The localization-status report is regenerated: every cell that still says
# "Ready for handoff" moves on to "In Translation". Because the new status
# text is shorter, the Status column(s) that held it are re-sized to fit
# the new (narrower) contents, just like the reporting job does when it
# rebuilds this workbook from the latest handoff data.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target "fit" width (in the stored character-width units Excel writes to
# <col width=".."/>) that the regenerated report ends up with for any
# column whose only content was the status text above. Excel's COM
# ColumnWidth setter only resolves to whole-pixel steps, so we solve for
# the ColumnWidth (in points-ish "characters") that rounds to the pixel
# step nearest this target and drive the object model with that number,
# rather than poking the worksheet XML directly.
$fitWidth = 13.4101845877511
$fitColumnWidth = ($fitWidth * 6 - 5) / 6

function Update-StatusCells {
    param($range)

    foreach ($cell in $range.Cells) {
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overviewStatusRange = $overview.Range("E2:F3")
Update-StatusCells $overviewStatusRange
$overviewStatusRange.ColumnWidth = $fitColumnWidth

# --- Per-locale detail sheets: Status column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $statusRange = $ws.Range("C2:C3")
    Update-StatusCells $statusRange
    $statusRange.ColumnWidth = $fitColumnWidth
}
